{"js": "// Add \"TEL\u00c9FONO: {{ tel }}\" and \"EMAIL: {{ email }}\" paragraphs right\n// after the \"DIRECCI\u00d3N FISCAL: {{ cxAddress }}\" paragraph, inside the\n// same table cell (client-info cell of the 2nd table in the document).\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\n// The client info (CLIENTE / RAZ\u00d3N SOCIAL / DIRECCI\u00d3N FISCAL) lives in\n// the first cell of the second table in the document (index 1).\nconst infoTable = tables.items[1];\nconst infoCell = infoTable.getCell(0, 0);\ninfoCell.body.paragraphs.load(\"items/text\");\nawait context.sync();\n\n// Locate the \"DIRECCI\u00d3N FISCAL\" paragraph (the last paragraph in the cell).\nlet addressParagraph = null;\nconst paras = infoCell.body.paragraphs.items;\nfor (let i = 0; i < paras.length; i++) {\n  if (paras[i].text.indexOf(\"DIRECCI\u00d3N FISCAL\") !== -1) {\n    addressParagraph = paras[i];\n    break;\n  }\n}\nif (!addressParagraph) {\n  addressParagraph = paras[paras.length - 1];\n}\n\nconst insertionRange = addressParagraph.getRange(\"End\");\n\n// Build the two new paragraphs as raw OOXML (wrapped in the Flat-OPC\n// \"pkg:package\" envelope that insertOoxml expects) so the exact\n// run-splitting / proofErr / language-tag markup is reproduced.\nconst ooxml = `<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>\n<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">\n  <pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\" pkg:padding=\"512\">\n    <pkg:xmlData>\n      <w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">\n        <w:body>\n          <w:p>\n            <w:pPr>\n              <w:pStyle w:val=\"TableParagraph\"/>\n              <w:spacing w:before=\"0\" w:line=\"210\" w:lineRule=\"exact\"/>\n              <w:jc w:val=\"left\"/>\n              <w:rPr>\n                <w:b/>\n                <w:sz w:val=\"18\"/>\n                <w:lang w:val=\"en-US\"/>\n              </w:rPr>\n            </w:pPr>\n            <w:r>\n              <w:rPr>\n                <w:b/>\n                <w:sz w:val=\"18\"/>\n              </w:rPr>\n              <w:t>TEL</w:t>\n            </w:r>\n            <w:r>\n              <w:rPr>\n                <w:b/>\n                <w:sz w:val=\"18\"/>\n                <w:lang w:val=\"es-419\"/>\n              </w:rPr>\n              <w:t xml:space=\"preserve\">\u00c9FONO: </w:t>\n            </w:r>\n            <w:proofErr w:type=\"gramStart\"/>\n            <w:r>\n              <w:rPr>\n                <w:b/>\n                <w:sz w:val=\"18\"/>\n                <w:lang w:val=\"en-US\"/>\n              </w:rPr>\n              <w:t xml:space=\"preserve\">{{ </w:t>\n            </w:r>\n            <w:proofErr w:type=\"spellStart\"/>\n            <w:r>\n              <w:rPr>\n                <w:b/>\n                <w:sz w:val=\"18\"/>\n                <w:lang w:val=\"en-US\"/>\n              </w:rPr>\n              <w:t>tel</w:t>\n            </w:r>\n            <w:proofErr w:type=\"spellEnd\"/>\n            <w:proofErr w:type=\"gramEnd\"/>\n            <w:r>\n              <w:rPr>\n                <w:b/>\n                <w:sz w:val=\"18\"/>\n                <w:lang w:val=\"en-US\"/>\n              </w:rPr>\n              <w:t xml:space=\"preserve\"> }}</w:t>\n            </w:r>\n          </w:p>\n          <w:p>\n            <w:pPr>\n              <w:pStyle w:val=\"TableParagraph\"/>\n              <w:spacing w:before=\"0\" w:line=\"210\" w:lineRule=\"exact\"/>\n              <w:jc w:val=\"left\"/>\n              <w:rPr>\n                <w:b/>\n                <w:sz w:val=\"18\"/>\n                <w:lang w:val=\"es-419\"/>\n              </w:rPr>\n            </w:pPr>\n            <w:r>\n              <w:rPr>\n                <w:b/>\n                <w:sz w:val=\"18\"/>\n                <w:lang w:val=\"es-419\"/>\n              </w:rPr>\n              <w:t xml:space=\"preserve\">EMAIL: </w:t>\n            </w:r>\n            <w:proofErr w:type=\"gramStart\"/>\n            <w:r>\n              <w:rPr>\n                <w:b/>\n                <w:sz w:val=\"18\"/>\n                <w:lang w:val=\"es-419\"/>\n              </w:rPr>\n              <w:t>{{ email</w:t>\n            </w:r>\n            <w:proofErr w:type=\"gramEnd\"/>\n            <w:r>\n              <w:rPr>\n                <w:b/>\n                <w:sz w:val=\"18\"/>\n                <w:lang w:val=\"es-419\"/>\n              </w:rPr>\n              <w:t xml:space=\"preserve\"> }}</w:t>\n            </w:r>\n          </w:p>\n          <w:sectPr/>\n        </w:body>\n      </w:document>\n    </pkg:xmlData>\n  </pkg:part>\n</pkg:package>`;\n\ninsertionRange.insertOoxml(ooxml, \"After\");\nawait context.sync();\n", "ps1": "# Add \"TEL\u00c9FONO: {{ tel }}\" and \"EMAIL: {{ email }}\" paragraphs right\n# after the \"DIRECCI\u00d3N FISCAL: {{ cxAddress }}\" paragraph, inside the\n# same table cell (client-info cell of the 2nd table in the document).\n#\n# NOTE: we intentionally never touch $d.Tables(...) before walking\n# $d.Paragraphs by index - doing so throws off this host's paragraph\n# index/anchor tracking. A plain text search over $d.Paragraphs is\n# sufficient (and unambiguous) to locate the target paragraph.\n\n$d = $word.ActiveDocument\n\n$target = $null\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    $para = $d.Paragraphs.Item($i)\n    if ($para.Range.Text.Contains(\"FISCAL:\")) {\n        $target = $para\n        break\n    }\n}\n\nif ($target -eq $null) {\n    throw \"Could not find the 'DIRECCION FISCAL' paragraph\"\n}\n\n# Collapsed range positioned right before the paragraph's end-of-\n# paragraph mark, so the inserted XML lands inside the same table\n# cell, right after the existing text.\n$insertPos = $target.Range.End - 1\n$rng = $d.Range($insertPos, $insertPos)\n\n$ooxml = @\"\n<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">\n  <pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\" pkg:padding=\"512\">\n    <pkg:xmlData>\n      <w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">\n        <w:body>\n          <w:p>\n            <w:pPr>\n              <w:pStyle w:val=\"TableParagraph\"/>\n              <w:spacing w:before=\"0\" w:line=\"210\" w:lineRule=\"exact\"/>\n              <w:jc w:val=\"left\"/>\n              <w:rPr>\n                <w:b/>\n                <w:sz w:val=\"18\"/>\n                <w:lang w:val=\"en-US\"/>\n              </w:rPr>\n            </w:pPr>\n            <w:r>\n              <w:rPr>\n                <w:b/>\n                <w:sz w:val=\"18\"/>\n              </w:rPr>\n              <w:t>TEL</w:t>\n            </w:r>\n            <w:r>\n              <w:rPr>\n                <w:b/>\n                <w:sz w:val=\"18\"/>\n                <w:lang w:val=\"es-419\"/>\n              </w:rPr>\n              <w:t xml:space=\"preserve\">\u00c9FONO: </w:t>\n            </w:r>\n            <w:proofErr w:type=\"gramStart\"/>\n            <w:r>\n              <w:rPr>\n                <w:b/>\n                <w:sz w:val=\"18\"/>\n                <w:lang w:val=\"en-US\"/>\n              </w:rPr>\n              <w:t xml:space=\"preserve\">{{ </w:t>\n            </w:r>\n            <w:proofErr w:type=\"spellStart\"/>\n            <w:r>\n              <w:rPr>\n                <w:b/>\n                <w:sz w:val=\"18\"/>\n                <w:lang w:val=\"en-US\"/>\n              </w:rPr>\n              <w:t>tel</w:t>\n            </w:r>\n            <w:proofErr w:type=\"spellEnd\"/>\n            <w:proofErr w:type=\"gramEnd\"/>\n            <w:r>\n              <w:rPr>\n                <w:b/>\n                <w:sz w:val=\"18\"/>\n                <w:lang w:val=\"en-US\"/>\n              </w:rPr>\n              <w:t xml:space=\"preserve\"> }}</w:t>\n            </w:r>\n          </w:p>\n          <w:p>\n            <w:pPr>\n              <w:pStyle w:val=\"TableParagraph\"/>\n              <w:spacing w:before=\"0\" w:line=\"210\" w:lineRule=\"exact\"/>\n              <w:jc w:val=\"left\"/>\n              <w:rPr>\n                <w:b/>\n                <w:sz w:val=\"18\"/>\n                <w:lang w:val=\"es-419\"/>\n              </w:rPr>\n            </w:pPr>\n            <w:r>\n              <w:rPr>\n                <w:b/>\n                <w:sz w:val=\"18\"/>\n                <w:lang w:val=\"es-419\"/>\n              </w:rPr>\n              <w:t xml:space=\"preserve\">EMAIL: </w:t>\n            </w:r>\n            <w:proofErr w:type=\"gramStart\"/>\n            <w:r>\n              <w:rPr>\n                <w:b/>\n                <w:sz w:val=\"18\"/>\n                <w:lang w:val=\"es-419\"/>\n              </w:rPr>\n              <w:t>{{ email</w:t>\n            </w:r>\n            <w:proofErr w:type=\"gramEnd\"/>\n            <w:r>\n              <w:rPr>\n                <w:b/>\n                <w:sz w:val=\"18\"/>\n                <w:lang w:val=\"es-419\"/>\n              </w:rPr>\n              <w:t xml:space=\"preserve\"> }}</w:t>\n            </w:r>\n          </w:p>\n          <w:sectPr/>\n        </w:body>\n      </w:document>\n    </pkg:xmlData>\n  </pkg:part>\n</pkg:package>\n\"@\n\n$rng.InsertXML($ooxml)\n"}
